$d = $word.ActiveDocument

# --- Change 1: "E* = r2 ..." -> the leading "E" was retyped as lowercase
#     "e", which splits the run into "e" + "* = r2 ...".
#     (Word splits a run whenever a bookmark boundary falls inside it, so we
#     use a scratch bookmark purely to force the split, then remove it.)
$rng1 = $d.Content
$null = $rng1.Find.Execute("E* = r2", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
$startE = $rng1.Start
$d.Range($startE, $startE + 1).Text = "e"
$splitRange1 = $d.Range($startE + 1, $startE + 1)
$d.Bookmarks.Add("_GoBack", $splitRange1)
$d.Bookmarks.Item("_GoBack").Delete()

# --- Change 2: "2 = particle two velocity" -> split into "2 = par" +
#     "ticle two velocity", with the "_GoBack" bookmark left sitting at the
#     split point. Adding "_GoBack" here also moves it away from its old
#     spot right after "Kd" (Word only ever keeps a single "_GoBack").
$rng2 = $d.Content
$null = $rng2.Find.Execute("2 = particle two velocity", $true, $false, `
                            $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint2 = $rng2.Start + "2 = par".Length
$splitRange2 = $d.Range($splitPoint2, $splitPoint2)
$d.Bookmarks.Add("_GoBack", $splitRange2)
